$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record progress updates for the MATLAB install task
$ws.Range("I6").Value = "Project runs, Serial comm issues"
$ws.Range("I5").Value = "Resubmitted, was upgraded to 2008b"

# Move the active selection to I7
$ws.Range("I7").Select()
